# Improve the function for the graphs representing cars that have not been
# fired: collapse the per-car columns (D:J) down to a single summary value
# in column C for each metric row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns D through J entirely (dimension shrinks from A1:J6 to A1:C6)
$ws.Range("D1:J6").EntireColumn.Delete()

# Update the remaining column C summary values for rows 2-6
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 55.9
$ws.Range("C4").Value = 24.23
$ws.Range("C5").Value = 1508.5
$ws.Range("C6").Value = 10
